$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "60.740.59"
$ws.Range("E2").Value = "  -3.71%  "

Set-TextValue $ws.Range("D3") "2.993.16"
$ws.Range("E3").Value = "  -5.84%  "

Set-TextValue $ws.Range("D5") "561.19"
$ws.Range("E5").Value = "  -5.15%  "

Set-TextValue $ws.Range("D6") "127.37"
$ws.Range("E6").Value = "  -6.46%  "

Set-TextValue $ws.Range("D8") "2.990.38"
$ws.Range("E8").Value = "  -5.86%  "

$ws.Range("E9").Value = "  -3.36%  "

$ws.Range("E10").Value = "  -6.25%  "

$ws.Range("E11").Value = "  -1.88%  "

Set-TextValue $ws.Range("D12") "0.436"
$ws.Range("E12").Value = "  -4.47%  "

$ws.Range("E13").Value = "  -6.77%  "

$ws.Range("E14").Value = "  -5.96%  "

Set-TextValue $ws.Range("D15") "0.119"
$ws.Range("E15").Value = "  +0.38%  "

Set-TextValue $ws.Range("D16") "3.485.79"
$ws.Range("E16").Value = "  -5.85%  "

Set-TextValue $ws.Range("D17") "60.868.75"
$ws.Range("E17").Value = "  -3.42%  "

Set-TextValue $ws.Range("D18") "2.993.33"
$ws.Range("E18").Value = "  -5.71%  "

$ws.Range("E19").Value = "  -6.93%  "

Set-TextValue $ws.Range("D20") "433.66"
$ws.Range("E20").Value = "  -5.99%  "

Set-TextValue $ws.Range("D21") "13.10"
$ws.Range("E21").Value = "  -5.88%  "

$ws.Range("E22").Value = "  -7.38%  "

$ws.Range("E23").Value = "  -7.15%  "

Set-TextValue $ws.Range("D24") "12.83"
$ws.Range("E24").Value = "  -4.54%  "

Set-TextValue $ws.Range("D25") "78.75"
$ws.Range("E25").Value = "  -5.64%  "

$ws.Range("E26").Value = "  +0.00%  "

Set-TextValue $ws.Range("D27") "1.00"
$ws.Range("E27").Value = "  -0.08%  "

$ws.Range("E28").Value = "  -7.52%  "

Set-TextValue $ws.Range("D29") "7.16"
$ws.Range("E29").Value = "  -8.34%  "

$ws.Range("E30").Value = "  -7.91%  "

Set-TextValue $ws.Range("D31") "25.41"
$ws.Range("E31").Value = "  -7.21%  "

$ws.Range("E32").Value = "  -11.39%  "

Set-TextValue $ws.Range("D33") "0.0931"
$ws.Range("E33").Value = "  -10.12%  "

$ws.Range("E34").Value = "  -4.91%  "

Set-TextValue $ws.Range("D35") "0.952"
$ws.Range("E35").Value = "  -8.50%  "

$ws.Range("E36").Value = "  -4.47%  "

Set-TextValue $ws.Range("D37") "49.90"
$ws.Range("E37").Value = "  -2.41%  "

$ws.Range("E38").Value = "  -6.61%  "

$ws.Range("E39").Value = "  -8.17%  "

Set-TextValue $ws.Range("D40") "7.76"
$ws.Range("E40").Value = "  -4.29%  "

$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue $ws.Range("D41") "0.107"
$ws.Range("E41").Value = "  -4.98%  "

$ws.Range("B42").Value = "Bittensor"
$ws.Range("C42").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue $ws.Range("D42") "373.90"
$ws.Range("E42").Value = "  -7.92%  "

Set-TextValue $ws.Range("D43") "2.681.35"
$ws.Range("E43").Value = "  -3.80%  "

$ws.Range("E44").Value = "  -9.08%  "

$ws.Range("E45").Value = "  +0.07%  "

Set-TextValue $ws.Range("D46") "120.57"
$ws.Range("E46").Value = "  -2.10%  "

$ws.Range("B47").Value = "TheGraph"
$ws.Range("C47").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
Set-TextValue $ws.Range("D47") "0.234"
$ws.Range("E47").Value = "  -7.37%  "

$ws.Range("B48").Value = "Arweave"
$ws.Range("C48").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
Set-TextValue $ws.Range("D48") "33.88"
$ws.Range("E48").Value = "  -2.45%  "

Set-TextValue $ws.Range("D49") "1.97"
$ws.Range("E49").Value = "  -8.29%  "

$ws.Range("E50").Value = "  -4.90%  "

Set-TextValue $ws.Range("D51") "23.33"
$ws.Range("E51").Value = "  -9.16%  "
